$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Autor column header already exists (R1 = "Autor") ---
# The validation/download flow now stamps every data row with the
# author who validated it, and stores "Fecha" as plain text instead of
# a date serial (the date widget posts its value back as a string).

# Columns whose values look numeric/date-like need to be forced to
# Text *before* the value is written, otherwise Excel auto-converts
# "12/02/2021" -> a date serial and "354135" -> a number. Resetting the
# style back to Normal afterwards keeps the cell format default.
function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Fecha (F) for existing rows 2-4: replace date serial with text ---
Set-TextValue $ws.Range("F2") "12/02/2021"
Set-TextValue $ws.Range("F3") "12/02/2021"
Set-TextValue $ws.Range("F4") "12/02/2021"

# --- Autor (R) for existing rows 2-4 ---
$ws.Range("R2").Value = "judith"
$ws.Range("R3").Value = "judith"

# Row 4 also got re-validated (Validado flips 0 -> 1) as part of the
# validate-button fix.
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = "judith"

# --- New rows pulled in by the download, rows 5-7 ---
$newRows = @(
    @{ Row=5; A=4; B=354135; C="120912lpñl"; D="354135";  E="poderes de los pobres"; F="12/02/2021"; G=1354;   H="industrial"; I=354354;  J="no se";  K="mexico"; L="uabc"; M=1; N=2; O=1; P=2; Q=0; R="judith" },
    @{ Row=6; A=5; B=65465;  C="peridos";    D="454654";  E="algo ninteres";         F="12/02/2021"; G=543645; H="sector2";   I=654654;  J="Poder"; K="EUA";    L="uabc"; M=2; N=2; O=1; P=1; Q=0; R="judith" },
    @{ Row=7; A=6; B=65432;  C="periodo555"; D="541354";  E="poderesss";             F="12/02/2021"; G=54354;  H="sector23";  I=321354;  J="mexico"; K="mexico"; L="uabc"; M=1; N=2; O=2; P=1; Q=0; R="judith" }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    Set-TextValue $ws.Cells.Item($row, 4) $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    Set-TextValue $ws.Cells.Item($row, 6) $r.F
    $ws.Cells.Item($row, 7).Value = $r.G
    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
}
